$wb = $excel.ActiveWorkbook

# "entidade" sheet (sheet1): update a couple of cells in row 2
$wsEntidade = $wb.Worksheets.Item("entidade")
$wsEntidade.Range("D2").Value = "Contribuinte ICMS"
$wsEntidade.Range("E2").Value = 1231231231234

# Selection/active-sheet bookkeeping: "entidade" becomes the active tab with
# E8 selected (was "setor" tab selected before, with S15 active on "entidade").
$wsEntidade.Activate() | Out-Null
$wsEntidade.Range("E8").Select() | Out-Null
